$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the two runs "What is the most common ticket class and type" + "?"
#    into a single run "What is the most common ticket class and type?"
#    (occurs twice in the document). A Find/Replace (ReplaceAll) over the
#    identical visible text collapses the run boundary at each match.
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("What is the most common ticket class and type?", $true, $false, $false, $false, $false, $true, 1, $false, "What is the most common ticket class and type?", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Append the new "EXCEL" analysis section at the end of the document.
#    The last paragraph in the document is an empty paragraph right before
#    the sectPr; replacing its Range with the OOXML fragment below inserts
#    11 blank paragraphs, the new EXCEL section, and re-creates the final
#    blank paragraph that Word always keeps before the section break.
# ---------------------------------------------------------------------------
$payloadXml = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"/>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"/>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"/>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"/>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"/>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"/>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"/>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"/>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"/>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"/>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"/>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>EXCEL</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:r>
        <w:t xml:space="preserve">Replaced </w:t>
      </w:r>
      <w:r>
        <w:t>Null Values in [</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos Narrow" w:eastAsia="Times New Roman" w:hAnsi="Aptos Narrow" w:cs="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="12501A"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>Actual</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos Narrow" w:eastAsia="Times New Roman" w:hAnsi="Aptos Narrow" w:cs="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="12501A"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> Arrival </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos Narrow" w:eastAsia="Times New Roman" w:hAnsi="Aptos Narrow" w:cs="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="12501A"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>Time</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos Narrow" w:eastAsia="Times New Roman" w:hAnsi="Aptos Narrow" w:cs="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="12501A"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">  ]</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos Narrow" w:eastAsia="Times New Roman" w:hAnsi="Aptos Narrow" w:cs="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="12501A"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>T</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">o </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">“ </w:t>
      </w:r>
      <w:r>
        <w:t>Cancelled - No Arrival</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> “ if the [ journey status ] is Cancelled  ( as a R</w:t>
      </w:r>
      <w:r>
        <w:t>eference</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> )</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:b/>
          <w:bCs/>
          <w:rtl/>
          <w:lang w:bidi="ar-EG"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:r>
        <w:t xml:space="preserve">Replaced </w:t>
      </w:r>
      <w:r>
        <w:t>Blank Values</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>in [</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos Narrow" w:eastAsia="Times New Roman" w:hAnsi="Aptos Narrow" w:cs="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="12501A"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> Reason</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos Narrow" w:eastAsia="Times New Roman" w:hAnsi="Aptos Narrow" w:cs="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="12501A"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> for </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos Narrow" w:eastAsia="Times New Roman" w:hAnsi="Aptos Narrow" w:cs="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="12501A"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>Delay</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos Narrow" w:eastAsia="Times New Roman" w:hAnsi="Aptos Narrow" w:cs="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="12501A"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos Narrow" w:eastAsia="Times New Roman" w:hAnsi="Aptos Narrow" w:cs="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="12501A"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> ]</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos Narrow" w:eastAsia="Times New Roman" w:hAnsi="Aptos Narrow" w:cs="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="12501A"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos Narrow" w:eastAsia="Times New Roman" w:hAnsi="Aptos Narrow" w:cs="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="12501A"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">to </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> “ No Delay “</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> ( </w:t>
      </w:r>
      <w:r>
        <w:t>as a R</w:t>
      </w:r>
      <w:r>
        <w:t>eference</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> )</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"/>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:r>
        <w:t xml:space="preserve">Added New Column [ </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>Delay ]</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Calculate the delay time for each ticket</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Aptos Narrow" w:eastAsia="Times New Roman" w:hAnsi="Aptos Narrow" w:cs="Times New Roman" w:hint="cs"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="12501A"/>
          <w:rtl/>
          <w:lang w:bidi="ar-EG"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:r>
        <w:t xml:space="preserve">Added New Column </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>[</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>Month</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Of </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>Purchase</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>]</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>to simplify further analysis</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"/>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:r>
        <w:t xml:space="preserve">Added New Column </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t xml:space="preserve">[ </w:t>
      </w:r>
      <w:r>
        <w:t>Day</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Of </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>Purchase</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> ]</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>to simplify further analysis</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"/>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:r>
        <w:t xml:space="preserve">Added New Column </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>[ Day</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Of </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>Journey</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> ]</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> to simplify further analysis</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"/>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Aptos Narrow" w:eastAsia="Times New Roman" w:hAnsi="Aptos Narrow" w:cs="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="12501A"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"/>

'@

$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertXML($payloadXml) | Out-Null
